$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")
$ws.Activate()

# Update the requirement text in B4 to split "자전거 상세 정보 조회" / "자전거 리스트 조회" use cases.
$ws.Range("B4").Value = "관리자는 자전거 리스트 조회 화면에서 원하는 자전거 항목을 선택해서 상세내용(자전거 ID, 자전거 제품명, 유형(일반/전기), 소속 대여소, 상태(사용 가능/수리 중))화면을 볼 수 있다."

# Move the active selection to B4 (matching the saved view state).
$ws.Range("B4").Select()
